$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(4652, 4652, 4664, 4706, 4936, 4964, 4976, 4976, 4976, 4976, 4976, 4980, 4996, 5066)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
